# Sivananda App.docx edit script
#
# Summary of the change (from commit message):
#   "fix issue with 60 minute class not showing the right asanas
#    fixed numberChooser being slow
#    increase numberChooser button width"
#
# Concretely, in the Todo list this means striking-through (marking done)
# a batch of previously-unstruck bullet items, appending a short note to
# one of them, rewording another, and moving the (hidden) "_GoBack" last
# edit position bookmark down to sit right after the newly-edited
# "Play 60,90,120 class" bullet (it previously sat right after the
# "Good to Have:" heading further down the document).

$d = $word.ActiveDocument

function Strike-Paragraph($para) {
    # Applying StrikeThrough to the whole paragraph range (which includes
    # the trailing paragraph mark) makes Word stamp <w:strike/> onto both
    # the paragraph-mark run properties (w:pPr/w:rPr) and the run(s) of
    # visible text -- matching how the rest of this list is already
    # formatted.
    $para.Range.Font.StrikeThrough = 1
}

function Find-ParagraphByText($containsText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($containsText)) {
            return $p
        }
    }
    return $null
}

# 1. "Some asanas don't have links" -- strike it, then append a new run
#    ", its ok" (also struck) as a separate run.
$p1 = Find-ParagraphByText("Some asanas don")
Strike-Paragraph $p1
$body1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$insertPos = $body1.End
$body1.InsertAfter(", its ok")
$newRun1 = $d.Range($insertPos, $body1.End)
# Toggling Bold on/off (net no-op) forces the engine to split this off
# into its own <w:r> instead of merging it into the previous run's text.
$newRun1.Font.Bold = 1
$newRun1.Font.Bold = 0

# 2. "Change tab bar icons to something useful"
Strike-Paragraph (Find-ParagraphByText("Change tab bar icons"))

# 3. "Number chooser is slow to respond"
Strike-Paragraph (Find-ParagraphByText("Number chooser is slow"))

# 4. "show prayer text when prayer is playing instead of image of asana"
Strike-Paragraph (Find-ParagraphByText("show prayer text when prayer is playing"))

# 5. "image of asana was removed... ignore this for now"
Strike-Paragraph (Find-ParagraphByText("image of asana was removed"))

# 6. "sound stops on phone off"
Strike-Paragraph (Find-ParagraphByText("sound stops on phone off"))

# 7. "phone is kept awake during playback"
Strike-Paragraph (Find-ParagraphByText("phone is kept awake during playback"))

# 8. "Play 60,90,120 class" -- strike it, and relocate the hidden
#    "_GoBack" bookmark here (right after this bullet's text).
$p8 = Find-ParagraphByText("Play 60,90,120")
Strike-Paragraph $p8

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$markPos = $p8.Range.End - 1
$marker = $d.Range($markPos, $markPos)
$marker.InsertAfter("_TMPMARK_")
$markerRange = $d.Range($markPos, $markPos + 9)
$d.Bookmarks.Add("_GoBack", $markerRange)
$cleanup = $d.Range($markPos, $markPos + 9)
$cleanup.Text = ""

# 9. "Create a screen"
Strike-Paragraph (Find-ParagraphByText("Create a screen"))

# 10. "Either combine audio files or play with a list" -- reworded to
#     "Play in a list", and struck through.
$p10 = Find-ParagraphByText("Either combine audio files")
$body10 = $d.Range($p10.Range.Start, $p10.Range.End - 1)
$body10.Text = "Play in a list"
Strike-Paragraph $p10
